# Generate Report for Handback
# Re-orders the four tracked files so the two that have been handed back
# (546360ff..., d1c15b0a...) move to the top of each sheet with their new
# "Handed back: in sync with en-US" status, and refreshes the handoff/handback
# file + datetime columns to match.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2:D5").Hyperlinks.Delete()

$ov.Range("A2").Value = "546360ff-766e-4363-a017-b55ff5251884.md"
$ov.Range("B2").Value = $handedBack
$ov.Range("C2").Value = $handedBack
$ov.Range("D2").Value = "2016-03-24 04:19:31"

$ov.Range("A3").Value = "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md"
$ov.Range("B3").Value = $handedBack
$ov.Range("C3").Value = $handedBack
$ov.Range("D3").Value = "2016-03-24 04:19:31"

$ov.Range("A4").Value = "e6b03b36-bd7f-4063-afe8-246553bc847e.md"
$ov.Range("B4").Value = "In Translation"
$ov.Range("C4").Value = "In Translation"
$ov.Range("D4").Value = "2016-03-24 04:18:00"

$ov.Range("A5").Value = "cee44cbc-cf37-454b-aa53-c86ea1678fd0.md"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"
$ov.Range("D5").Value = "2016-03-24 04:19:31"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/546360ff-766e-4363-a017-b55ff5251884.md", "", "", "546360ff-766e-4363-a017-b55ff5251884.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/d1c15b0a-4d52-4643-93fb-6fca5d885c58.md", "", "", "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb2ca0e3632d2300aeddb117afaad6b8601545e/e2e/e6b03b36-bd7f-4063-afe8-246553bc847e.md", "", "", "e6b03b36-bd7f-4063-afe8-246553bc847e.md")
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7e21c860ce9aa1f0db5e607c7a83bbd13d4c8baa/e2e/cee44cbc-cf37-454b-aa53-c86ea1678fd0.md", "", "", "cee44cbc-cf37-454b-aa53-c86ea1678fd0.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2:L5").Hyperlinks.Delete()

# wipe the F/G helper columns first; only rows 2 & 3 keep data in the new layout
$zh.Range("F2:G5").ClearContents()

$zh.Range("A2").Value = "546360ff-766e-4363-a017-b55ff5251884.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $handedBack
$zh.Range("D2").Value = "546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-24 04:19:27"
$zh.Range("F2").Value = "546360ff-766e-4363-a017-b55ff5251884.md"
$zh.Range("G2").Value = "546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-24 04:19:50"
$zh.Range("J2").Value = "Include"

$zh.Range("A3").Value = "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = $handedBack
$zh.Range("D3").Value = "d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-24 04:19:27"
$zh.Range("F3").Value = "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md"
$zh.Range("G3").Value = "d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-24 04:19:50"
$zh.Range("J3").Value = "Include"

$zh.Range("A4").Value = "e6b03b36-bd7f-4063-afe8-246553bc847e.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "In Translation"
$zh.Range("D4").Value = "e6b03b36-bd7f-4063-afe8-246553bc847e.11f3e460b8ea2de1ee338f993ee8d53c6a1e1ff2.zh-cn.xlf"
$zh.Range("E4").Value = "2016-03-24 04:17:56"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("J4").Value = "Include"

$zh.Range("A5").Value = "cee44cbc-cf37-454b-aa53-c86ea1678fd0.md"
$zh.Range("B5").Value = ".md"
$zh.Range("C5").Value = "Ready for handoff"
$zh.Range("D5").Value = "cee44cbc-cf37-454b-aa53-c86ea1678fd0.abb9b124f1a14c1d0fc2182472ced6242c0d9bec.zh-cn.xlf"
$zh.Range("E5").Value = "2016-03-24 04:19:27"
$zh.Range("H5").Value = "0001-01-01 00:00:00"
$zh.Range("J5").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/546360ff-766e-4363-a017-b55ff5251884.md", "", "", "546360ff-766e-4363-a017-b55ff5251884.md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0796967ed1e74ec460c0dbf990862683af687609/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.zh-cn.xlf", "", "", "546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/546360ff-766e-4363-a017-b55ff5251884.md", "", "", "546360ff-766e-4363-a017-b55ff5251884.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0796967ed1e74ec460c0dbf990862683af687609/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.zh-cn.xlf", "", "", "546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/d1c15b0a-4d52-4643-93fb-6fca5d885c58.md", "", "", "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0796967ed1e74ec460c0dbf990862683af687609/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.zh-cn.xlf", "", "", "d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/d1c15b0a-4d52-4643-93fb-6fca5d885c58.md", "", "", "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md")
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0796967ed1e74ec460c0dbf990862683af687609/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.zh-cn.xlf", "", "", "d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb2ca0e3632d2300aeddb117afaad6b8601545e/e2e/e6b03b36-bd7f-4063-afe8-246553bc847e.md", "", "", "e6b03b36-bd7f-4063-afe8-246553bc847e.md")
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90340a5b685bf34446a32d0cc8e370fb7a9ccdc4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e6b03b36-bd7f-4063-afe8-246553bc847e.11f3e460b8ea2de1ee338f993ee8d53c6a1e1ff2.zh-cn.xlf", "", "", "e6b03b36-bd7f-4063-afe8-246553bc847e.11f3e460b8ea2de1ee338f993ee8d53c6a1e1ff2.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7e21c860ce9aa1f0db5e607c7a83bbd13d4c8baa/e2e/cee44cbc-cf37-454b-aa53-c86ea1678fd0.md", "", "", "cee44cbc-cf37-454b-aa53-c86ea1678fd0.md")
$zh.Hyperlinks.Add($zh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0796967ed1e74ec460c0dbf990862683af687609/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/cee44cbc-cf37-454b-aa53-c86ea1678fd0.abb9b124f1a14c1d0fc2182472ced6242c0d9bec.zh-cn.xlf", "", "", "cee44cbc-cf37-454b-aa53-c86ea1678fd0.abb9b124f1a14c1d0fc2182472ced6242c0d9bec.zh-cn.xlf")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2:L5").Hyperlinks.Delete()

$de.Range("F2:G5").ClearContents()

$de.Range("A2").Value = "546360ff-766e-4363-a017-b55ff5251884.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $handedBack
$de.Range("D2").Value = "546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.de-de.xlf"
$de.Range("E2").Value = "2016-03-24 04:19:31"
$de.Range("F2").Value = "546360ff-766e-4363-a017-b55ff5251884.md"
$de.Range("G2").Value = "546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.de-de.xlf"
$de.Range("H2").Value = "2016-03-24 04:19:56"
$de.Range("J2").Value = "Include"

$de.Range("A3").Value = "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = $handedBack
$de.Range("D3").Value = "d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.de-de.xlf"
$de.Range("E3").Value = "2016-03-24 04:19:31"
$de.Range("F3").Value = "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md"
$de.Range("G3").Value = "d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.de-de.xlf"
$de.Range("H3").Value = "2016-03-24 04:19:56"
$de.Range("J3").Value = "Include"

$de.Range("A4").Value = "e6b03b36-bd7f-4063-afe8-246553bc847e.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "In Translation"
$de.Range("D4").Value = "e6b03b36-bd7f-4063-afe8-246553bc847e.11f3e460b8ea2de1ee338f993ee8d53c6a1e1ff2.de-de.xlf"
$de.Range("E4").Value = "2016-03-24 04:18:00"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("J4").Value = "Include"

$de.Range("A5").Value = "cee44cbc-cf37-454b-aa53-c86ea1678fd0.md"
$de.Range("B5").Value = ".md"
$de.Range("C5").Value = "Ready for handoff"
$de.Range("D5").Value = "cee44cbc-cf37-454b-aa53-c86ea1678fd0.abb9b124f1a14c1d0fc2182472ced6242c0d9bec.de-de.xlf"
$de.Range("E5").Value = "2016-03-24 04:19:31"
$de.Range("H5").Value = "0001-01-01 00:00:00"
$de.Range("J5").Value = "Include"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/546360ff-766e-4363-a017-b55ff5251884.md", "", "", "546360ff-766e-4363-a017-b55ff5251884.md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/935dc114ef6625eaa26ad83c1db2a1fdbbf91a03/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.de-de.xlf", "", "", "546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/546360ff-766e-4363-a017-b55ff5251884.md", "", "", "546360ff-766e-4363-a017-b55ff5251884.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/935dc114ef6625eaa26ad83c1db2a1fdbbf91a03/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.de-de.xlf", "", "", "546360ff-766e-4363-a017-b55ff5251884.71a3efbd0cf4dd33689f6a1d3b93cac448d40a5b.de-de.xlf")

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/d1c15b0a-4d52-4643-93fb-6fca5d885c58.md", "", "", "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/935dc114ef6625eaa26ad83c1db2a1fdbbf91a03/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.de-de.xlf", "", "", "d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5bbd8bf322124fe14b88bfd811926b62a3f76b/e2e/d1c15b0a-4d52-4643-93fb-6fca5d885c58.md", "", "", "d1c15b0a-4d52-4643-93fb-6fca5d885c58.md")
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/935dc114ef6625eaa26ad83c1db2a1fdbbf91a03/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.de-de.xlf", "", "", "d1c15b0a-4d52-4643-93fb-6fca5d885c58.6d036b0efa4f22ec930816952fb489e8dccb5ca5.de-de.xlf")

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb2ca0e3632d2300aeddb117afaad6b8601545e/e2e/e6b03b36-bd7f-4063-afe8-246553bc847e.md", "", "", "e6b03b36-bd7f-4063-afe8-246553bc847e.md")
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2ba6e7adc0b5a7fdc4437f412d070f034ed797b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e6b03b36-bd7f-4063-afe8-246553bc847e.11f3e460b8ea2de1ee338f993ee8d53c6a1e1ff2.de-de.xlf", "", "", "e6b03b36-bd7f-4063-afe8-246553bc847e.11f3e460b8ea2de1ee338f993ee8d53c6a1e1ff2.de-de.xlf")

$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7e21c860ce9aa1f0db5e607c7a83bbd13d4c8baa/e2e/cee44cbc-cf37-454b-aa53-c86ea1678fd0.md", "", "", "cee44cbc-cf37-454b-aa53-c86ea1678fd0.md")
$de.Hyperlinks.Add($de.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/935dc114ef6625eaa26ad83c1db2a1fdbbf91a03/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/cee44cbc-cf37-454b-aa53-c86ea1678fd0.abb9b124f1a14c1d0fc2182472ced6242c0d9bec.de-de.xlf", "", "", "cee44cbc-cf37-454b-aa53-c86ea1678fd0.abb9b124f1a14c1d0fc2182472ced6242c0d9bec.de-de.xlf")

"Report regenerated for handback"
